# HW3 document revision: "Further explanation in the doc file"
# Applies the textual / structural edits described by the target diff
# using Word COM-interop (Find/Execute + Range surgery).

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output ("NOT FOUND: " + $old)
    }
}

# 1. "Chen " + "Naveh" (two runs, spell-checked) -> single run "Chen Naveh"
Replace-Text "Chen Naveh" "Chen Naveh"

# 2. Window-size paragraph: rework the Taylor-series clause
Replace-Text "the construction of the optical flow formulas won't be correct especially when using Taylor series" "the Taylor series used in the construction of the optical flow formulas will be less accurate and therefore the approximation will be damaged"

# 3. K (distance between frames) paragraph - full rewrite
Replace-Text "When increasing K (the distance between frames) and there is a fast motion we will miss those motion therefore for scenarios we want to segment background there is a chance we will mark moving objects as background." "When increasing K (the distance between frames) we sample the video in a lower frequency which can lead to loosing fast motions. When using all frames we get better accuracy but also more noise as we detect minor changes in the background. Also using all frames will have effect on the performance of the algorithm."

# 4. Merge the two "scale" bullets into a single, rewritten bullet
$pLarger = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("If we use larger")) {
        $pLarger = $p
        break
    }
}
if ($pLarger -ne $null) {
    $pNext = $pLarger.Next()
    $mergedRange = $d.Range($pLarger.Range.Start, $pNext.Range.End)
    $mergedRange.Text = "The scale parameter affects the pixel resolution of the image. In small scales the image pixels appear bigger and the image is more blurred, this allow us to detect changes in large items in the image. As the scale get closer to 1 we get the original image which contain more pixels and more fine details. This allow us to detect minor changes and is more susceptible to noise. The right order to use will be small scale first and large scale at the end. This allows us to get major movements first and then to refine the changes be using higher resolution."
    # Re-find the (now single) paragraph and remove the now-empty trailing paragraph left behind
    $pLarger2 = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.StartsWith("The scale parameter affects")) {
            $pLarger2 = $p
            break
        }
    }
    if ($pLarger2 -ne $null) {
        $pAfter = $pLarger2.Next()
        if ($pAfter -ne $null -and $pAfter.Range.Text.Trim().Length -eq 0) {
            $pAfter.Range.Delete()
        }
    }
} else {
    Write-Output "NOT FOUND: If we use larger scale bullet"
}

# 5. "A region OF works better..." (children/slide) bullet - full rewrite
Replace-Text "A region OF works better is at the top of the slide region. Most of the video children are sitting on top of the slide and when they move the algorithm recognize the slide behind them we obviously the slide show be marked as background and not foreground. The reason it is like that in the median detector is because most of the video the pixels didn't change and then when the top slide was revealed the detector thought it is a foreground" "A region OF works better is at the top of the slide region. Most of the video the children are sitting on top of the slide, when they move the algorithm recognize the slide behind them. We obviously would think of the slide as the background and the children as foreground. This is an example for the median detector disadvantage as it estimates the background as the pixels median values (in our case the children) and the foreground as changes which appear for a short part of the video (the slide)."

# 6. Question 16 answer (assumptions of optical flow) - full rewrite
Replace-Text "The assumption we take into consideration when using optical flow are that the motion in the scene should be small and the intensity of the pixels should be consistency. Moreover, we will receive better results when each pixel is moving in the same directions as its neighbors" "The assumption we make when using magnitude thresholding is that the motion around each pixel is in the same direction (for the Taylor series) and intensity of the pixels should be consistent (we assume the objects/pixels do not change this also true for the illumination in the image). "

# 7. "For instance ... move fast ..." -> "moves fast"
Replace-Text "For instance assume two objects in the 3D world, one which is further away from the camera and move fast and the other which is closer and moves slower. both can be detected with the same optical flow." "For instance assume two objects in the 3D world, one which is further away from the camera and moves fast and the other which is closer and moves slower. both can be detected with the same optical flow."

# 8. "assume two objects one which has vector component..." paragraph - rewrite
Replace-Text "assume two objects one which has vector component away from the camera (away from center of projection) and the other which doesn't. both can be detected with the same optical flow even though they move directionally different in the 3D world" "assume two objects one which has vector component away from the camera (away from center of projection) and the other which doesn't. Both can be detected with the same optical flow (same vector direction in this case) even though they move directionally different in the 3D world."

# 9. Question 18 answer (planar scene) - full rewrite
Replace-Text "In order to recognize if a scene is planar we can film the scene by moving slowly the camera and calculating the optical flow the video. If orientation of the optical flow is equal for most of the scene then we can assume the scene is planar" "In order to determine if a scene is planar we can apply the OF algorithm to the video. If orientation of the optical flow is equal for most of the scene then we can assume the scene is planar. The intuition is that in a video of a planar surface while the camera is moving, we expect most of the image to move in the opposite direction of the camera. So if the camera is moving to the left, we expect to see the patterns in the video moving to the right."

# 10. "axsis" typo fix + extended explanation
Replace-Text "The expected orientation of the optical flow will be similar in all pixels and in the X axsis" "The expected orientation of the optical flow will be similar in all pixels and in the X axis. This is due to the fact that the scene is static, therefore all movements are the consequences of the moving camera. In this case the camera moves in the X axis and so all the scene should be moving in the X direction. No movement in the Y direction is expected."

Write-Output "done"
